$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the Date column (BF) values: the original text "2-16-2013-14" was a
# mangled representation of the game date; correct it to ISO "2014-02-16"
# for every data row (rows 2 through 31), keeping the values as plain text
# (not auto-converted to an Excel date serial number).
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    if ($cell.Value2 -eq "2-16-2013-14") {
        $cell.NumberFormat = "@"
        $cell.Value = "2014-02-16"
        $cell.Style = "Normal"
    }
}
